$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C4").Value = "http://twitter.github.com/bootstrap/javascript.html#tabs"
$ws.Range("B4").Value = "active menu"
